# one_unique_organisations export: the org_id (column B) values for rows
# 2-24 got reshuffled (a pairwise swap between the top half and bottom half
# of the table), while the org_name (column C) stayed attached to its row -
# except for the two rows whose underlying organisation identity changed,
# which also got a corrected org_name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row pairs whose org_id (column B) values trade places.
$pairs = @(
    @(2, 14),
    @(3, 24),
    @(4, 23),
    @(5, 22),
    @(6, 21),
    @(7, 20),
    @(8, 19),
    @(9, 18),
    @(10, 17),
    @(11, 16),
    @(12, 15)
)

foreach ($pair in $pairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $cell1 = $ws.Cells.Item($r1, 2)
    $cell2 = $ws.Cells.Item($r2, 2)

    $b1 = $cell1.Value2
    $b2 = $cell2.Value2

    # org_id is stored as text (e.g. "11374"), not a number - force the
    # Text number format before assigning so the swapped ids keep their
    # original text type instead of being reinterpreted as numerics.
    $cell1.NumberFormat = "@"
    $cell1.Value2 = $b2

    $cell2.NumberFormat = "@"
    $cell2.Value2 = $b1
}

# Row 11 now carries org_id 570, whose name is "ИПНГ РАН".
$ws.Cells.Item(11, 3).Value2 = "ИПНГ РАН"

# Row 16 now carries org_id 4493, whose name is the English variant.
$ws.Cells.Item(16, 3).Value2 = "Oil and Gas Research Institute|Russian Academy of Sciences (OGRI RAS)"
